$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Model_Home_win" (column C) and "Model_home_win_probability" (column D)
# predictions for rows 2-49, reflecting the refreshed model run.
$data = @(
    @(2, 1, 0.5099221467971802),
    @(3, 0, 0.2572522759437561),
    @(4, 0, 0.4490349590778351),
    @(5, 0, 0.4267487823963165),
    @(6, 0, 0.2212837189435959),
    @(7, 1, 0.6441329717636108),
    @(8, 0, 0.3838265836238861),
    @(9, 1, 0.516806960105896),
    @(10, 0, 0.2984212636947632),
    @(11, 1, 0.7743140459060669),
    @(12, 1, 0.7111673355102539),
    @(13, 0, 0.3285875022411346),
    @(14, 1, 0.9595517516136169),
    @(15, 1, 0.6709197759628296),
    @(16, 0, 0.4069570004940033),
    @(17, 0, 0.4987024068832397),
    @(18, 0, 0.4332131445407867),
    @(19, 0, 0.409697026014328),
    @(20, 0, 0.2965400516986847),
    @(21, 0, 0.362677663564682),
    @(22, 1, 0.5705922842025757),
    @(23, 0, 0.1554374098777771),
    @(24, 1, 0.6417803764343262),
    @(25, 0, 0.4861744344234467),
    @(26, 1, 0.6391391754150391),
    @(27, 0, 0.4958766400814056),
    @(28, 0, 0.05065657570958138),
    @(29, 1, 0.5099936127662659),
    @(30, 0, 0.4267343580722809),
    @(31, 1, 0.7769376635551453),
    @(32, 0, 0.4972838163375854),
    @(33, 0, 0.2526662945747375),
    @(34, 0, 0.0960831344127655),
    @(35, 1, 0.9502853155136108),
    @(36, 0, 0.3619717657566071),
    @(37, 1, 0.5929375290870667),
    @(38, 1, 0.5340604782104492),
    @(39, 0, 0.2213118076324463),
    @(40, 0, 0.4904601871967316),
    @(41, 1, 0.8556072115898132),
    @(42, 1, 0.6698856353759766),
    @(43, 1, 0.9412673711776733),
    @(44, 1, 0.7091012597084045),
    @(45, 0, 0.3101401925086975),
    @(46, 1, 0.7661099433898926),
    @(47, 0, 0.4192559421062469),
    @(48, 1, 0.6347099542617798),
    @(49, 1, 0.922538697719574)
)

foreach ($row in $data) {
    $r = $row[0]
    $cVal = $row[1]
    $dVal = $row[2]
    $ws.Cells.Item($r, 3).Value = $cVal
    $ws.Cells.Item($r, 4).Value = $dVal
}
